$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Sheet "Logs": append a new row (42) with a new test-mail log entry
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A42").Value = "Zou jij dit even op kunnen pakken?"
$logs.Range("B42").Value = "mailmind.test@zohomail.eu"
$logs.Range("C42").Value = "Testmail #1: Zou jij dit even op kunnen pakken?"
$logs.Range("D42").Value = "Overig"
$logs.Range("E42").Value = "Geachte afzender,`nDank u voor uw bericht. Kunt u meer details geven over wat u precies wilt dat we oppakken? Zo kunnen we u beter van dienst zijn.`nMet vriendelijke groet,`n[Naam van het bedrijf]"
$logs.Range("F42").Value = "2025-08-05 19:21:21"
$logs.Range("G42").Value = "Ja"
$logs.Range("H42").Value = "Nee"
$logs.Range("I42").Value = "Ja"
$logs.Range("J42").Value = "Nee"

# Extend the conditional formatting ranges on "Logs" so they cover the new row
$logsFcD = $logs.Range("D2:D41").FormatConditions
for ($i = 1; $i -le $logsFcD.Count; $i++) {
    $logsFcD.Item($i).ModifyAppliesToRange($logs.Range("D2:D42"))
}

$logsFcG = $logs.Range("G2:G41").FormatConditions
for ($i = 1; $i -le $logsFcG.Count; $i++) {
    $logsFcG.Item($i).ModifyAppliesToRange($logs.Range("G2:G42"))
}

$logsFcH = $logs.Range("H2:H41").FormatConditions
for ($i = 1; $i -le $logsFcH.Count; $i++) {
    $logsFcH.Item($i).ModifyAppliesToRange($logs.Range("H2:H42"))
}

$logsFcI = $logs.Range("I2:I41").FormatConditions
for ($i = 1; $i -le $logsFcI.Count; $i++) {
    $logsFcI.Item($i).ModifyAppliesToRange($logs.Range("I2:I42"))
}

$logsFcJ = $logs.Range("J2:J41").FormatConditions
for ($i = 1; $i -le $logsFcJ.Count; $i++) {
    $logsFcJ.Item($i).ModifyAppliesToRange($logs.Range("J2:J42"))
}

# ---------------------------------------------------------------------
# 2) Sheet "Dashboard": swap category order of rows 3 & 4, add row 10
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A3").Value = "Inkoop / Bestellingen"
$dash.Range("A4").Value = "Klantenservice / Contact"

$dash.Range("A10").Value = "Overig"
$dash.Range("B10").Value = 1

# ---------------------------------------------------------------------
# 3) Update the chart series on "Dashboard" to include the new row 10
# ---------------------------------------------------------------------
$chartObj = $dash.ChartObjects(1)
$chart = $chartObj.Chart
$ser = $chart.SeriesCollection(1)
$ser.Formula = "=SERIES(Dashboard!`$B`$1,Dashboard!`$A`$2:`$A`$10,Dashboard!`$B`$2:`$B`$10,1)"

$wb.Save()
